$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 146.17053
$ws.Range("H2").Value = 438.51159
$ws.Range("I2").Value = 0.4047435297111188
$ws.Range("J2").Value = 0.4047435297111188
$ws.Range("M2").Value = 159.4836373333333
$ws.Range("N2").Value = 478.450912
$ws.Range("O2").Value = 0.2983285084902258
$ws.Range("P2").Value = 0.2983285084902258
$ws.Range("Q2").Value = 23311.80779534112
$ws.Range("R2").Value = 209806.2701580701
$ws.Range("S2").Value = 0.1207465335397875
$ws.Range("T2").Value = 0.1207465335397875
$ws.Range("G3").Value = 146.17053
$ws.Range("H3").Value = 438.51159
$ws.Range("I3").Value = 0.4047435297111188
$ws.Range("J3").Value = 0.4047435297111188
$ws.Range("O3").Value = 0.3227862111630279
$ws.Range("P3").Value = 0.3227862111630279
$ws.Range("Q3").Value = 25222.96696249341
$ws.Range("R3").Value = 227006.7026624407
$ws.Range("S3").Value = 0.1306456304482024
$ws.Range("T3").Value = 0.1306456304482024
$ws.Range("G4").Value = 146.17053
$ws.Range("H4").Value = 438.51159
$ws.Range("I4").Value = 0.4047435297111188
$ws.Range("J4").Value = 0.4047435297111188
$ws.Range("M4").Value = 74.38770566666666
$ws.Range("N4").Value = 223.163117
$ws.Range("O4").Value = 0.1391489036280481
$ws.Range("P4").Value = 0.1391489036280482
$ws.Range("Q4").Value = 10873.29036278067
$ws.Range("R4").Value = 97859.61326502603
$ws.Range("S4").Value = 0.0563196184098485
$ws.Range("T4").Value = 0.05631961840984851
$ws.Range("G5").Value = 146.17053
$ws.Range("H5").Value = 438.51159
$ws.Range("I5").Value = 0.4047435297111188
$ws.Range("J5").Value = 0.4047435297111188
$ws.Range("M5").Value = 58.41461433333333
$ws.Range("N5").Value = 175.243843
$ws.Range("O5").Value = 0.1092697975759847
$ws.Range("P5").Value = 0.1092697975759848
$ws.Range("Q5").Value = 8538.495136848929
$ws.Range("R5").Value = 76846.45623164036
$ws.Range("S5").Value = 0.04422624356172351
$ws.Range("T5").Value = 0.04422624356172352
$ws.Range("G6").Value = 146.17053
$ws.Range("H6").Value = 438.51159
$ws.Range("I6").Value = 0.4047435297111188
$ws.Range("J6").Value = 0.4047435297111188
$ws.Range("M6").Value = 69.746216
$ws.Range("N6").Value = 209.238648
$ws.Range("O6").Value = 0.1304665791427133
$ws.Range("P6").Value = 0.1304665791427133
$ws.Range("Q6").Value = 10194.84135821448
$ws.Range("R6").Value = 91753.57222393031
$ws.Range("S6").Value = 0.05280550375155681
$ws.Range("T6").Value = 0.05280550375155683
$ws.Range("H7").Value = 632.3552549999999
$ws.Range("I7").Value = 0.5836600531814327
$ws.Range("J7").Value = 0.5836600531814327
$ws.Range("M7").Value = 159.4836373333333
$ws.Range("N7").Value = 478.450912
$ws.Range("O7").Value = 0.2983285084902258
$ws.Range("P7").Value = 0.2983285084902258
$ws.Range("Q7").Value = 33616.77205141584
$ws.Range("R7").Value = 302550.9484627426
$ws.Range("S7").Value = 0.1741224331309427
$ws.Range("T7").Value = 0.1741224331309427
$ws.Range("H8").Value = 632.3552549999999
$ws.Range("I8").Value = 0.5836600531814327
$ws.Range("J8").Value = 0.5836600531814327
$ws.Range("O8").Value = 0.3227862111630279
$ws.Range("P8").Value = 0.3227862111630279
$ws.Range("R8").Value = 327354.8171185551
$ws.Range("S8").Value = 0.188397417173646
$ws.Range("T8").Value = 0.188397417173646
$ws.Range("H9").Value = 632.3552549999999
$ws.Range("I9").Value = 0.5836600531814327
$ws.Range("J9").Value = 0.5836600531814327
$ws.Range("M9").Value = 74.38770566666666
$ws.Range("N9").Value = 223.163117
$ws.Range("O9").Value = 0.1391489036280481
$ws.Range("P9").Value = 0.1391489036280482
$ws.Range("Q9").Value = 15679.81886190331
$ws.Range("R9").Value = 141118.3697571298
$ws.Range("S9").Value = 0.08121565649168462
$ws.Range("T9").Value = 0.08121565649168463
$ws.Range("H10").Value = 632.3552549999999
$ws.Range("I10").Value = 0.5836600531814327
$ws.Range("J10").Value = 0.5836600531814327
$ws.Range("M10").Value = 58.41461433333333
$ws.Range("N10").Value = 175.243843
$ws.Range("O10").Value = 0.1092697975759847
$ws.Range("P10").Value = 0.1092697975759848
$ws.Range("Q10").Value = 12312.92944749388
$ws.Range("R10").Value = 110816.365027445
$ws.Range("S10").Value = 0.06377641586432364
$ws.Range("T10").Value = 0.06377641586432364
$ws.Range("H11").Value = 632.3552549999999
$ws.Range("I11").Value = 0.5836600531814327
$ws.Range("J11").Value = 0.5836600531814327
$ws.Range("M11").Value = 69.746216
$ws.Range("N11").Value = 209.238648
$ws.Range("O11").Value = 0.1304665791427133
$ws.Range("P11").Value = 0.1304665791427133
$ws.Range("Q11").Value = 14701.46206798836
$ws.Range("R11").Value = 132313.1586118952
$ws.Range("S11").Value = 0.07614813052083565
$ws.Range("T11").Value = 0.07614813052083566
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.5494936666666667
$ws.Range("H12").Value = 1.648481
$ws.Range("I12").Value = 0.00152153793381314
$ws.Range("J12").Value = 0.00152153793381314
$ws.Range("M12").Value = 159.4836373333333
$ws.Range("N12").Value = 478.450912
$ws.Range("O12").Value = 0.2983285084902258
$ws.Range("P12").Value = 0.2983285084902258
$ws.Range("Q12").Value = 87.63524865163023
$ws.Range("R12").Value = 788.717237864672
$ws.Range("S12").Value = 0.000453918142405774
$ws.Range("T12").Value = 0.000453918142405774
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.5494936666666667
$ws.Range("H13").Value = 1.648481
$ws.Range("I13").Value = 0.00152153793381314
$ws.Range("J13").Value = 0.00152153793381314
$ws.Range("O13").Value = 0.3227862111630279
$ws.Range("P13").Value = 0.3227862111630279
$ws.Range("Q13").Value = 94.81980123101899
$ws.Range("R13").Value = 853.378211079171
$ws.Range("S13").Value = 0.0004911314647963655
$ws.Range("T13").Value = 0.0004911314647963655
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.5494936666666667
$ws.Range("H14").Value = 1.648481
$ws.Range("I14").Value = 0.00152153793381314
$ws.Range("J14").Value = 0.00152153793381314
$ws.Range("M14").Value = 74.38770566666666
$ws.Range("N14").Value = 223.163117
$ws.Range("O14").Value = 0.1391489036280481
$ws.Range("P14").Value = 0.1391489036280482
$ws.Range("Q14").Value = 40.87557314169744
$ws.Range("R14").Value = 367.880158275277
$ws.Range("S14").Value = 0.0002117203353185841
$ws.Range("T14").Value = 0.0002117203353185842
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.5494936666666667
$ws.Range("H15").Value = 1.648481
$ws.Range("I15").Value = 0.00152153793381314
$ws.Range("J15").Value = 0.00152153793381314
$ws.Range("M15").Value = 58.41461433333333
$ws.Range("N15").Value = 175.243843
$ws.Range("O15").Value = 0.1092697975759847
$ws.Range("P15").Value = 0.1092697975759848
$ws.Range("Q15").Value = 32.09846061694255
$ws.Range("R15").Value = 288.886145552483
$ws.Range("S15").Value = 0.0001662581420319439
$ws.Range("T15").Value = 0.0001662581420319439
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.5494936666666667
$ws.Range("H16").Value = 1.648481
$ws.Range("I16").Value = 0.00152153793381314
$ws.Range("J16").Value = 0.00152153793381314
$ws.Range("M16").Value = 69.746216
$ws.Range("N16").Value = 209.238648
$ws.Range("O16").Value = 0.1304665791427133
$ws.Range("P16").Value = 0.1304665791427133
$ws.Range("Q16").Value = 38.32510396596533
$ws.Range("R16").Value = 344.925935693688
$ws.Range("S16").Value = 0.0001985098492604726
$ws.Range("T16").Value = 0.0001985098492604726
$ws.Range("G17").Value = 3.410044
$ws.Range("H17").Value = 10.230132
$ws.Range("I17").Value = 0.009442349596941478
$ws.Range("J17").Value = 0.009442349596941478
$ws.Range("M17").Value = 159.4836373333333
$ws.Range("N17").Value = 478.450912
$ws.Range("O17").Value = 0.2983285084902258
$ws.Range("P17").Value = 0.2983285084902258
$ws.Range("Q17").Value = 543.8462205867094
$ws.Range("R17").Value = 4894.615985280384
$ws.Range("S17").Value = 0.002816922071898836
$ws.Range("T17").Value = 0.002816922071898836
$ws.Range("G18").Value = 3.410044
$ws.Range("H18").Value = 10.230132
$ws.Range("I18").Value = 0.009442349596941478
$ws.Range("J18").Value = 0.009442349596941478
$ws.Range("O18").Value = 0.3227862111630279
$ws.Range("P18").Value = 0.3227862111630279
$ws.Range("Q18").Value = 588.4320673438679
$ws.Range("R18").Value = 5295.888606094812
$ws.Range("S18").Value = 0.003047860250873483
$ws.Range("T18").Value = 0.003047860250873483
$ws.Range("G19").Value = 3.410044
$ws.Range("H19").Value = 10.230132
$ws.Range("I19").Value = 0.009442349596941478
$ws.Range("J19").Value = 0.009442349596941478
$ws.Range("M19").Value = 74.38770566666666
$ws.Range("N19").Value = 223.163117
$ws.Range("O19").Value = 0.1391489036280481
$ws.Range("P19").Value = 0.1391489036280482
$ws.Range("Q19").Value = 253.6653493823826
$ws.Range("R19").Value = 2282.988144441444
$ws.Range("S19").Value = 0.001313892594087149
$ws.Range("T19").Value = 0.001313892594087149
$ws.Range("G20").Value = 3.410044
$ws.Range("H20").Value = 10.230132
$ws.Range("I20").Value = 0.009442349596941478
$ws.Range("J20").Value = 0.009442349596941478
$ws.Range("M20").Value = 58.41461433333333
$ws.Range("N20").Value = 175.243843
$ws.Range("O20").Value = 0.1092697975759847
$ws.Range("P20").Value = 0.1092697975759848
$ws.Range("Q20").Value = 199.1964051196973
$ws.Range("R20").Value = 1792.767646077276
$ws.Range("S20").Value = 0.001031763629099476
$ws.Range("T20").Value = 0.001031763629099477
$ws.Range("G21").Value = 3.410044
$ws.Range("H21").Value = 10.230132
$ws.Range("I21").Value = 0.009442349596941478
$ws.Range("J21").Value = 0.009442349596941478
$ws.Range("M21").Value = 69.746216
$ws.Range("N21").Value = 209.238648
$ws.Range("O21").Value = 0.1304665791427133
$ws.Range("P21").Value = 0.1304665791427133
$ws.Range("Q21").Value = 237.837665393504
$ws.Range("R21").Value = 2140.538988541536
$ws.Range("S21").Value = 0.001231911050982532
$ws.Range("T21").Value = 0.001231911050982533
$ws.Range("E22").Value = 3
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 0.228434
$ws.Range("H22").Value = 0.6853020000000001
$ws.Range("I22").Value = 0.0006325295766939459
$ws.Range("J22").Value = 0.0006325295766939459
$ws.Range("M22").Value = 159.4836373333333
$ws.Range("N22").Value = 478.450912
$ws.Range("O22").Value = 0.2983285084902258
$ws.Range("P22").Value = 0.2983285084902258
$ws.Range("Q22").Value = 36.43148521060267
$ws.Range("R22").Value = 327.883366895424
$ws.Range("S22").Value = 0.0001887016051910588
$ws.Range("T22").Value = 0.0001887016051910588
$ws.Range("E23").Value = 3
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 0.228434
$ws.Range("H23").Value = 0.6853020000000001
$ws.Range("I23").Value = 0.0006325295766939459
$ws.Range("J23").Value = 0.0006325295766939459
$ws.Range("O23").Value = 0.3227862111630279
$ws.Range("P23").Value = 0.3227862111630279
$ws.Range("Q23").Value = 39.418227703698
$ws.Range("R23").Value = 354.764049333282
$ws.Range("S23").Value = 0.0002041718255095927
$ws.Range("T23").Value = 0.0002041718255095927
$ws.Range("E24").Value = 3
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = 0.228434
$ws.Range("H24").Value = 0.6853020000000001
$ws.Range("I24").Value = 0.0006325295766939459
$ws.Range("J24").Value = 0.0006325295766939459
$ws.Range("M24").Value = 74.38770566666666
$ws.Range("N24").Value = 223.163117
$ws.Range("O24").Value = 0.1391489036280481
$ws.Range("P24").Value = 0.1391489036280482
$ws.Range("Q24").Value = 16.99268115625933
$ws.Range("R24").Value = 152.934130406334
$ws.Range("S24").Value = 0.00008801579710927595
$ws.Range("T24").Value = 0.00008801579710927598
$ws.Range("E25").Value = 3
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = 0.228434
$ws.Range("H25").Value = 0.6853020000000001
$ws.Range("I25").Value = 0.0006325295766939459
$ws.Range("J25").Value = 0.0006325295766939459
$ws.Range("M25").Value = 58.41461433333333
$ws.Range("N25").Value = 175.243843
$ws.Range("O25").Value = 0.1092697975759847
$ws.Range("P25").Value = 0.1092697975759848
$ws.Range("Q25").Value = 13.34388401062067
$ws.Range("R25").Value = 120.094956095586
$ws.Range("S25").Value = 0.00006911637880617079
$ws.Range("T25").Value = 0.00006911637880617079
$ws.Range("E26").Value = 3
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 0.228434
$ws.Range("H26").Value = 0.6853020000000001
$ws.Range("I26").Value = 0.0006325295766939459
$ws.Range("J26").Value = 0.0006325295766939459
$ws.Range("M26").Value = 69.746216
$ws.Range("N26").Value = 209.238648
$ws.Range("O26").Value = 0.1304665791427133
$ws.Range("P26").Value = 0.1304665791427133
$ws.Range("Q26").Value = 15.932407105744
$ws.Range("R26").Value = 143.391663951696
$ws.Range("S26").Value = 0.00008252397007784765
$ws.Range("T26").Value = 0.00008252397007784766
